# Refresh the "cryptos" price list (Price / Volume(1h) columns, plus the
# Kaspa <-> EthereumClassic row swap) to match the latest scrape.
#
# Several Price values are plain-looking decimals (e.g. "384.70", "0.0320").
# Excel's normal text-entry parsing would coerce those to numbers and drop
# the significant trailing zero, so for those cells we force the cell to
# Text format first (NumberFormat = "@") and then assign the literal string,
# exactly as typing '384.70 into the cell would. Values that already contain
# two dots (e.g. "51.711.25") are never number-like and don't need this.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.711.25'
$ws.Range("E2").Value = '  +1.17%  '
$ws.Range("D3").Value = '3.065.87'
$ws.Range("E3").Value = '  +3.59%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '384.70'
$ws.Range("E5").Value = '  +1.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.39'
$ws.Range("E6").Value = '  +1.20%  '
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.588'
$ws.Range("E9").Value = '  -0.73%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.21'
$ws.Range("E10").Value = '  +2.23%  '
$ws.Range("E11").Value = '  +0.46%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0866'
$ws.Range("E12").Value = '  +1.09%  '
$ws.Range("D13").Value = '3.545.50'
$ws.Range("E13").Value = '  +3.56%  '
$ws.Range("E14").Value = '  +2.34%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.77'
$ws.Range("E15").Value = '  -0.86%  '
$ws.Range("D16").Value = '3.083.69'
$ws.Range("E16").Value = '  +3.92%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.982'
$ws.Range("E17").Value = '  -1.43%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.55'
$ws.Range("E18").Value = '  -6.48%  '
$ws.Range("D19").Value = '51.771.46'
$ws.Range("E19").Value = '  +1.15%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.14'
$ws.Range("E20").Value = '  +0.10%  '
$ws.Range("E21").Value = '  +1.25%  '
$ws.Range("E22").Value = '  +0.36%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.19'
$ws.Range("E23").Value = '  -0.33%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '269.60'
$ws.Range("E24").Value = '  +1.08%  '
$ws.Range("E25").Value = '  -3.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.49'
$ws.Range("E26").Value = '  +8.31%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '27.26'
$ws.Range("E27").Value = '  +5.40%  '
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.173'
$ws.Range("E28").Value = '  +5.18%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.30'
$ws.Range("E29").Value = '  +1.38%  '
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("E31").Value = '  -1.47%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '10.33'
$ws.Range("E32").Value = '  +0.49%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.52'
$ws.Range("E33").Value = '  +0.52%  '
$ws.Range("E34").Value = '  +0.68%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '50.48'
$ws.Range("E35").Value = '  -1.40%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0445'
$ws.Range("E36").Value = '  +2.14%  '
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.41'
$ws.Range("E38").Value = '  +5.78%  '
$ws.Range("E39").Value = '  +3.99%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.287'
$ws.Range("E40").Value = '  +5.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.89'
$ws.Range("E41").Value = '  +3.00%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '128.71'
$ws.Range("E42").Value = '  +3.21%  '
$ws.Range("E43").Value = '  +0.07%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.54'
$ws.Range("E44").Value = '  +1.82%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.72'
$ws.Range("E45").Value = '  +5.00%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.22'
$ws.Range("E46").Value = '  +3.62%  '
$ws.Range("E47").Value = '  +7.70%  '
$ws.Range("E48").Value = '  +3.42%  '
$ws.Range("D49").Value = '2.048.38'
$ws.Range("E49").Value = '  +0.46%  '
$ws.Range("D50").Value = '3.368.47'
$ws.Range("E50").Value = '  +3.65%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0320'
$ws.Range("E51").Value = '  -0.04%  '
